$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 114000.3
$ws.Range("I43").Value = 4333.8335
$ws.Range("K43").Value = 4333.8335
$ws.Range("M43").Value = -4264.8335
# Row 57
$ws.Range("H57").Value = 30270
$ws.Range("J57").Value = 30270
$ws.Range("L57").Value = 90810
$ws.Range("N57").Value = -91808
# Row 98
$ws.Range("H98").Value = 2245.9
$ws.Range("I98").Value = 2063.625
$ws.Range("J98").Value = 2975
$ws.Range("K98").Value = 2063.625
$ws.Range("L98").Value = 2975
$ws.Range("M98").Value = -565.625
$ws.Range("N98").Value = -5971
# Row 116
$ws.Range("H116").Value = 82140.484
$ws.Range("I116").Value = 127630.18
$ws.Range("J116").Value = 4808
$ws.Range("K116").Value = 127630.18
$ws.Range("L116").Value = 4808
$ws.Range("M116").Value = -124188.18
$ws.Range("N116").Value = -11692
# Row 122
$ws.Range("H122").Value = 2245.9
$ws.Range("I122").Value = 2063.625
$ws.Range("J122").Value = 2975
$ws.Range("K122").Value = 6190.875
$ws.Range("L122").Value = 8925
$ws.Range("M122").Value = -3740.875
$ws.Range("N122").Value = -13825
# Row 132
$ws.Range("H132").Value = 3598.889
$ws.Range("I132").Value = 1940.4884
$ws.Range("K132").Value = 5821.4652
$ws.Range("M132").Value = -3291.4652
# Row 138
$ws.Range("H138").Value = 2751.6165
$ws.Range("I138").Value = 2318.1875
$ws.Range("J138").Value = 2873.2808
$ws.Range("K138").Value = 6954.5625
$ws.Range("L138").Value = 8619.8424
$ws.Range("M138").Value = -1814.5625
$ws.Range("N138").Value = -18899.8424

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1776
$ws.Range("I45").Value = 1594.8572
$ws.Range("J45").Value = 1957.1428
$ws.Range("K45").Value = 1594.8572
$ws.Range("L45").Value = 1957.1428
$ws.Range("M45").Value = -1217.8572
$ws.Range("N45").Value = -2711.1428
# Row 63
$ws.Range("H63").Value = 8660.799999999999
$ws.Range("I63").Value = 10829.714
$ws.Range("J63").Value = 3600
$ws.Range("K63").Value = 10829.714
$ws.Range("L63").Value = 3600
$ws.Range("M63").Value = -10143.714
$ws.Range("N63").Value = -4972
# Row 66
$ws.Range("H66").Value = 8660.799999999999
$ws.Range("I66").Value = 10829.714
$ws.Range("J66").Value = 3600
$ws.Range("K66").Value = 54148.57
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = -50716.57
$ws.Range("N66").Value = -24864
# Row 110
$ws.Range("H110").Value = 2341.5715
$ws.Range("I110").Value = 2315.1667
$ws.Range("J110").Value = 2500
$ws.Range("K110").Value = 2315.1667
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = -270.1667000000002
$ws.Range("N110").Value = -6590
# Row 132
$ws.Range("H132").Value = 2632.627
$ws.Range("I132").Value = 2460.7234
$ws.Range("J132").Value = 3036.6
$ws.Range("K132").Value = 7382.1702
$ws.Range("L132").Value = 9109.799999999999
$ws.Range("M132").Value = -4852.1702
$ws.Range("N132").Value = -14169.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
# Row 94
$ws.Range("H94").Value = 1774.091
$ws.Range("I94").Value = 1640.875
$ws.Range("K94").Value = 1640.875
$ws.Range("M94").Value = -1189.875
# Row 105
$ws.Range("H105").Value = 2349.8906
$ws.Range("I105").Value = 2088.5557
$ws.Range("K105").Value = 2088.5557
$ws.Range("M105").Value = -341.5556999999999
# Row 141
$ws.Range("H141").Value = 47773.332
$ws.Range("J141").Value = 47773.332
$ws.Range("L141").Value = 47773.332
$ws.Range("N141").Value = -58133.332

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 2617.5881
$ws.Range("I58").Value = 2031.4736
$ws.Range("J58").Value = 3360
$ws.Range("K58").Value = 2031.4736
$ws.Range("L58").Value = 3360
$ws.Range("M58").Value = -1828.4736
$ws.Range("N58").Value = -3766
# Row 99
$ws.Range("H99").Value = 45385.824
$ws.Range("I99").Value = 85050.336
$ws.Range("J99").Value = 2115.4546
$ws.Range("K99").Value = 85050.336
$ws.Range("L99").Value = 2115.4546
$ws.Range("M99").Value = -83552.336
$ws.Range("N99").Value = -5111.4546
# Row 126
$ws.Range("H126").Value = 45385.824
$ws.Range("I126").Value = 85050.336
$ws.Range("J126").Value = 2115.4546
$ws.Range("K126").Value = 255151.008
$ws.Range("L126").Value = 6346.3638
$ws.Range("M126").Value = -252681.008
$ws.Range("N126").Value = -11286.3638
# Row 134
$ws.Range("H134").Value = 11906621
$ws.Range("I134").Value = 18519866
$ws.Range("J134").Value = 2778.4666
$ws.Range("K134").Value = 55559598
$ws.Range("L134").Value = 8335.399800000001
$ws.Range("M134").Value = -55557063
$ws.Range("N134").Value = -13405.3998
# Row 136
$ws.Range("H136").Value = 2617.5881
$ws.Range("I136").Value = 2031.4736
$ws.Range("J136").Value = 3360
$ws.Range("K136").Value = 6094.4208
$ws.Range("L136").Value = 10080
$ws.Range("M136").Value = -3544.4208
$ws.Range("N136").Value = -15180

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 738.8095
$ws.Range("I122").Value = 411.83334
$ws.Range("J122").Value = 1174.7778
$ws.Range("K122").Value = 3706.50006
$ws.Range("L122").Value = 10573.0002
$ws.Range("M122").Value = -1256.50006
$ws.Range("N122").Value = -15473.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 109
$ws.Range("H109").Value = 34500
$ws.Range("J109").Value = 34500
$ws.Range("L109").Value = 34500
$ws.Range("N109").Value = -36580
# Row 132
$ws.Range("H132").Value = 2953.9077
$ws.Range("I132").Value = 2571.84
$ws.Range("J132").Value = 4227.467
$ws.Range("K132").Value = 7715.52
$ws.Range("L132").Value = 12682.401
$ws.Range("M132").Value = -5185.52
$ws.Range("N132").Value = -17742.401

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2481.6924
$ws.Range("I82").Value = 2439.087
$ws.Range("J82").Value = 2542.9375
$ws.Range("K82").Value = 2439.087
$ws.Range("L82").Value = 2542.9375
$ws.Range("M82").Value = -2078.087
$ws.Range("N82").Value = -3264.9375
# Row 85
$ws.Range("H85").Value = 2481.6924
$ws.Range("I85").Value = 2439.087
$ws.Range("J85").Value = 2542.9375
$ws.Range("K85").Value = 2439.087
$ws.Range("L85").Value = 2542.9375
$ws.Range("M85").Value = -1191.087
$ws.Range("N85").Value = -5038.9375
# Row 122
$ws.Range("H122").Value = 3174.875
$ws.Range("I122").Value = 3492
$ws.Range("J122").Value = 2857.75
$ws.Range("K122").Value = 10476
$ws.Range("L122").Value = 8573.25
$ws.Range("M122").Value = -8026
$ws.Range("N122").Value = -13473.25
# Row 132
$ws.Range("H132").Value = 7414448.5
$ws.Range("I132").Value = 2496.1428
$ws.Range("J132").Value = 19622370
$ws.Range("K132").Value = 7488.428400000001
$ws.Range("L132").Value = 58867110
$ws.Range("M132").Value = -4958.428400000001
$ws.Range("N132").Value = -58872170

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1592.3798
$ws.Range("I132").Value = 731.4151000000001
$ws.Range("J132").Value = 3347.423
$ws.Range("K132").Value = 2194.2453
$ws.Range("L132").Value = 10042.269
$ws.Range("M132").Value = 335.7547
$ws.Range("N132").Value = -15102.269
# Row 136
$ws.Range("H136").Value = 1954.2821
$ws.Range("I136").Value = 1103.963
$ws.Range("J136").Value = 3867.5
$ws.Range("K136").Value = 3311.889
$ws.Range("L136").Value = 11602.5
$ws.Range("M136").Value = -761.8890000000001
$ws.Range("N136").Value = -16702.5

# Special case: BSM row 9 -> clear N9 cell (removed in diff)
$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("N9").ClearContents()
